# Admin Copy to other outlet
# Re-create the view/formatting tweaks made to Sheet1:
#  - scroll the sheet so C6 becomes the top-left visible cell
#  - make the header row (row 1) taller (27pt, explicit custom height)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Keep the existing selection (A10) but move the viewport's top-left
# corner to C6, matching sheetView/@topLeftCell="C6" in the saved file.
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 3

# Header row gets an explicit 27pt height.
$ws.Rows.Item(1).RowHeight = 27
